$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.933
$ws.Range("C18").Value = -12.199
$ws.Range("A21").Value = -20.645
$ws.Range("A23").Value = -20.808
$ws.Range("B24").Value = 5.218999999999999
$ws.Range("A25").Value = -21.63
$ws.Range("B28").Value = 5.121
$ws.Range("B36").Value = 7.363
$ws.Range("B45").Value = 5.655999999999999
$ws.Range("B48").Value = 5.48
$ws.Range("B49").Value = 6.237
$ws.Range("C51").Value = -11.589
$ws.Range("B52").Value = 4.992
$ws.Range("A53").Value = -20.671
$ws.Range("B53").Value = 7.512
$ws.Range("B54").Value = 5.197
$ws.Range("C55").Value = -13.655
$ws.Range("A57").Value = -22.178
$ws.Range("A59").Value = -22.461
$ws.Range("C64").Value = -10.76
$ws.Range("A69").Value = -21.591
$ws.Range("B70").Value = 4.935
$ws.Range("A79").Value = -21.305
$ws.Range("C80").Value = -12.033
$ws.Range("A83").Value = -21.976
$ws.Range("B86").Value = 5.197
$ws.Range("B87").Value = 4.924000000000001
$ws.Range("C92").Value = -10.975
$ws.Range("A93").Value = -21.687
$ws.Range("C94").Value = -11.35
$ws.Range("C96").Value = -10.318
$ws.Range("B101").Value = 5.197000000000001
